$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '72.003.48'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.69%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '4.042.60'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.33%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '536.85'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.36%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '149.45'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.34%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '4.037.25'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.37%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.697'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.15%  '
$ws.Range("E9").Value = '  -0.06%  '
$ws.Range("E10").Value = '  -0.66%  '
$ws.Range("E11").Value = '  -0.84%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '53.38'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +11.44%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000330'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.35%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.90'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.58%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.684.78'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.22%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.036.76'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.17%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.30'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.13%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '20.68'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.34%  '
$ws.Range("E19").Value = '  +0.03%  '
$ws.Range("E20").Value = '  -0.88%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.031.70'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.58%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '436.83'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.91%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '98.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.84%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.50'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.08%  '
$ws.Range("E25").Value = '  +0.18%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '14.69'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.20%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.42'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +28.02%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.24'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.11%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.70'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.79%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.96'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.15%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '37.13'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.45%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.51'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +24.38%  '
$ws.Range("E33").Value = '  +3.14%  '
$ws.Range("B34").Value = 'Cosmos'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '13.57'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.41%  '
$ws.Range("B35").Value = 'InjectiveProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '49.91'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +17.96%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '679.25'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.42%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '66.66'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.66%  '
$ws.Range("E38").Value = '  +6.69%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0₃0895'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +6.15%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.46'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +7.50%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.149'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.54%  '
$ws.Range("E42").Value = '  -2.25%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.27'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +18.59%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.999'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("E45").Value = '  +0.02%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.999'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.09%  '
$ws.Range("E47").Value = '  -0.58%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.66'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.53%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.13'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.23%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.28'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.99%  '
$ws.Range("B51").Value = 'FLOKI'
$ws.Range("C51").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.000277'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.04%  '
